$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.923.43"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.292.07"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'505.72"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'129.47"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.312.77"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "'0.0979"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").Value = "'5.14"
$ws.Range("E12").Value = "  +8.56%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'23.72"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "2.702.48"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "54.941.17"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "2.312.25"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'10.50"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'4.18"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'310.88"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").Value = "'6.64"
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'60.31"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "'0.993"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'7.51"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").Value = "'172.11"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("D30").Value = "0.0₃0708"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "'0.925"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("D38").Value = "'3.88"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "'36.82"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'134.39"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'261.12"
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("D46").Value = "'0.0505"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").Value = "'16.49"
$ws.Range("E51").Value = "  +0.82%  "

# Clear quote-prefix formatting introduced by forcing numeric-looking text, to keep cell style index unchanged
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D51").ClearFormats()
